$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: OT values that were "Pendiente ADM" now have ICD numbers ---
$ws.Range("E27").Value = "ICD30794466"
$ws.Range("E30").Value = "ICD30794595"

# --- Helper: write a text value into a cell while keeping default (no explicit) style ---
function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

# --- Row 32 ---
Set-TextCell $ws.Range("A32") "7216"
Set-TextCell $ws.Range("B32") "9/15/2025"
Set-TextCell $ws.Range("C32") "NEWBERY, JORGE AV. 3870"
$ws.Range("D32").Value = 15
Set-TextCell $ws.Range("E32") "ICD30800947"
Set-TextCell $ws.Range("F32") "Optical Power"
Set-TextCell $ws.Range("G32") "Pendiente"
Set-TextCell $ws.Range("H32") "rienda cortada reparar o retirar"
$ws.Range("I32").Value = 1
Set-TextCell $ws.Range("J32") '{"direccionesNormalizadas": [{"altura": 3870, "cod_calle": 14019, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.450024", "y": "-34.586606"}, "direccion": "NEWBERY, JORGE AV. 3870, CABA", "nombre_calle": "NEWBERY, JORGE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K32").Value = -58.450024
$ws.Range("L32").Value = -34.586606
Set-TextCell $ws.Range("M32") "Colegiales"
Set-TextCell $ws.Range("N32") "Capital Norte"

# --- Row 33 ---
Set-TextCell $ws.Range("A33") "7217"
Set-TextCell $ws.Range("B33") "9/15/2025"
Set-TextCell $ws.Range("C33") "CARRANZA, ANGEL JUSTINIANO 1135"
$ws.Range("D33").Value = 15
Set-TextCell $ws.Range("E33") "ICD30801240"
Set-TextCell $ws.Range("F33") "Optical Power"
Set-TextCell $ws.Range("G33") "Pendiente"
Set-TextCell $ws.Range("H33") "Cable en panza"
$ws.Range("I33").Value = 1
Set-TextCell $ws.Range("J33") '{"direccionesNormalizadas": [{"altura": 1135, "cod_calle": 3074, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.443900", "y": "-34.587281"}, "direccion": "CARRANZA, ANGEL JUSTINIANO 1135, CABA", "nombre_calle": "CARRANZA, ANGEL JUSTINIANO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K33").Value = -58.4439
$ws.Range("L33").Value = -34.587281
Set-TextCell $ws.Range("M33") "Colegiales"
Set-TextCell $ws.Range("N33") "Capital Norte"

# --- Row 34 ---
Set-TextCell $ws.Range("A34") "7219"
Set-TextCell $ws.Range("B34") "9/15/2025"
Set-TextCell $ws.Range("C34") "CRAMER AV. 2064"
$ws.Range("D34").Value = 13
Set-TextCell $ws.Range("E34") "ICD30801472"
Set-TextCell $ws.Range("F34") "Optical Power"
Set-TextCell $ws.Range("G34") "Pendiente"
Set-TextCell $ws.Range("H34") "Rienda cortada cable suelto"
$ws.Range("I34").Value = 1
Set-TextCell $ws.Range("J34") '{"direccionesNormalizadas": [{"altura": 2064, "cod_calle": 3189, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.461063", "y": "-34.565190"}, "direccion": "CRAMER AV. 2064, CABA", "nombre_calle": "CRAMER AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K34").Value = -58.461063
$ws.Range("L34").Value = -34.56519
Set-TextCell $ws.Range("M34") "Colegiales"
Set-TextCell $ws.Range("N34") "Capital Norte"

# --- Row 35 ---
Set-TextCell $ws.Range("A35") "2099"
Set-TextCell $ws.Range("B35") "9/15/2025"
Set-TextCell $ws.Range("C35") "YRIGOYEN, HIPOLITO AV. 3818"
$ws.Range("D35").Value = 5
Set-TextCell $ws.Range("E35") "ICD30801756"
Set-TextCell $ws.Range("F35") "Optical Power"
Set-TextCell $ws.Range("G35") "Pendiente"
Set-TextCell $ws.Range("H35") "Cables a baja altura"
$ws.Range("I35").Value = 1
Set-TextCell $ws.Range("J35") '{"direccionesNormalizadas": [{"altura": 3818, "cod_calle": 26005, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.420687", "y": "-34.614694"}, "direccion": "YRIGOYEN, HIPOLITO AV. 3818, CABA", "nombre_calle": "YRIGOYEN, HIPOLITO AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K35").Value = -58.420687
$ws.Range("L35").Value = -34.614694
Set-TextCell $ws.Range("M35") "Almagro"
Set-TextCell $ws.Range("N35") "Capital Sur"
